$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1964285714285714
$ws.Range("C2").Value = 0.5586734693877551
$ws.Range("J2").Value = 0.02551020408163265
$ws.Range("P2").Value = 0.1352040816326531
$ws.Range("S2").Value = 0.08418367346938775
$ws.Range("B3").Value = 0.008928571428571428
$ws.Range("C3").Value = 0.008928571428571428
$ws.Range("J3").Value = 0.03125
$ws.Range("P3").Value = 0.7544642857142857
$ws.Range("S3").Value = 0.1964285714285714
$ws.Range("J4").Value = 0.03278688524590164
$ws.Range("P4").Value = 0.7049180327868853
$ws.Range("S4").Value = 0.2622950819672131
$ws.Range("B6").Value = 0.06198347107438017
$ws.Range("F6").Value = 0.04545454545454546
$ws.Range("J6").Value = 0.2809917355371901
$ws.Range("O6").Value = 0.01652892561983471
$ws.Range("Q6").Value = 0.2190082644628099
$ws.Range("R6").Value = 0.05785123966942149
$ws.Range("S6").Value = 0.3181818181818182
$ws.Range("B7").Value = 0.1170731707317073
$ws.Range("D7").Value = 0.04390243902439024
$ws.Range("F7").Value = 0.02926829268292683
$ws.Range("J7").Value = 0.1658536585365854
$ws.Range("O7").Value = 0.03902439024390244
$ws.Range("Q7").Value = 0.2146341463414634
$ws.Range("R7").Value = 0.08780487804878048
$ws.Range("S7").Value = 0.3024390243902439
$ws.Range("B8").Value = 0.1002004008016032
$ws.Range("D8").Value = 0.01803607214428858
$ws.Range("F8").Value = 0.05811623246492986
$ws.Range("J8").Value = 0.1082164328657315
$ws.Range("O8").Value = 0.03206412825651302
$ws.Range("Q8").Value = 0.2044088176352706
$ws.Range("R8").Value = 0.07214428857715431
$ws.Range("S8").Value = 0.406813627254509
$ws.Range("B9").Value = 0.1016042780748663
$ws.Range("D9").Value = 0.0267379679144385
$ws.Range("F9").Value = 0.09090909090909091
$ws.Range("J9").Value = 0.1122994652406417
$ws.Range("O9").Value = 0.03208556149732621
$ws.Range("Q9").Value = 0.1550802139037433
$ws.Range("R9").Value = 0.106951871657754
$ws.Range("S9").Value = 0.374331550802139
$ws.Range("B10").Value = 0.1325706594885599
$ws.Range("D10").Value = 0.02624495289367429
$ws.Range("F10").Value = 0.06527590847913863
$ws.Range("J10").Value = 0.1083445491251682
$ws.Range("O10").Value = 0.02826379542395693
$ws.Range("Q10").Value = 0.2388963660834455
$ws.Range("R10").Value = 0.07267833109017496
$ws.Range("S10").Value = 0.3277254374158816
$ws.Range("G11").Value = 0.1535947712418301
$ws.Range("J11").Value = 0.09477124183006536
$ws.Range("K11").Value = 0.2124183006535948
$ws.Range("L11").Value = 0.5261437908496732
$ws.Range("S11").Value = 0.0130718954248366
$ws.Range("G12").Value = 0.7298850574712644
$ws.Range("J12").Value = 0.2011494252873563
$ws.Range("K12").Value = 0.01149425287356322
$ws.Range("L12").Value = 0.04022988505747126
$ws.Range("S12").Value = 0.01724137931034483
$ws.Range("G13").Value = 0.603448275862069
$ws.Range("J13").Value = 0.3275862068965517
$ws.Range("S13").Value = 0.06896551724137931
$ws.Range("F15").Value = 0.02264150943396226
$ws.Range("H15").Value = 0.1320754716981132
$ws.Range("I15").Value = 0.04150943396226415
$ws.Range("J15").Value = 0.3169811320754717
$ws.Range("K15").Value = 0.06415094339622641
$ws.Range("M15").Value = 0.003773584905660377
$ws.Range("O15").Value = 0.04905660377358491
$ws.Range("S15").Value = 0.369811320754717
$ws.Range("F16").Value = 0.01953125
$ws.Range("H16").Value = 0.203125
$ws.Range("I16").Value = 0.0546875
$ws.Range("J16").Value = 0.453125
$ws.Range("K16").Value = 0.06640625
$ws.Range("M16").Value = 0.0234375
$ws.Range("O16").Value = 0.0390625
$ws.Range("S16").Value = 0.140625
$ws.Range("F17").Value = 0.02101576182136602
$ws.Range("H17").Value = 0.2101576182136602
$ws.Range("I17").Value = 0.0893169877408056
$ws.Range("J17").Value = 0.4430823117338004
$ws.Range("K17").Value = 0.08056042031523643
$ws.Range("M17").Value = 0.03502626970227671
$ws.Range("N17").Value = 0.001751313485113835
$ws.Range("O17").Value = 0.05078809106830123
$ws.Range("S17").Value = 0.06830122591943957
$ws.Range("F18").Value = 0.01538461538461539
$ws.Range("H18").Value = 0.1384615384615385
$ws.Range("I18").Value = 0.07692307692307693
$ws.Range("J18").Value = 0.4871794871794872
$ws.Range("K18").Value = 0.08205128205128205
$ws.Range("M18").Value = 0.05128205128205128
$ws.Range("O18").Value = 0.05641025641025641
$ws.Range("S18").Value = 0.09230769230769231
$ws.Range("F19").Value = 0.0178432893716059
$ws.Range("H19").Value = 0.2110162916989915
$ws.Range("I19").Value = 0.07525213343677269
$ws.Range("J19").Value = 0.4003103180760279
$ws.Range("K19").Value = 0.1031807602792863
$ws.Range("M19").Value = 0.01706749418153607
$ws.Range("N19").Value = 0.001551590380139643
$ws.Range("O19").Value = 0.07835531419705198
$ws.Range("S19").Value = 0.09542280837858805
